# ---------------------------------------------------------------------------
# 1) Re-style the three tables (slides 14, 15, 16) with the Table Styles
#    gallery entry that was picked in the Table Design ribbon.
# ---------------------------------------------------------------------------
$p = $ppt.ActivePresentation

$newTableStyleId = "{8DA8A0FB-8795-401A-B2E8-B330B24E3186}"

foreach ($slideIndex in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIndex)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shape = $slide.Shapes.Item($i)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newTableStyleId)
        }
    }
}

# ---------------------------------------------------------------------------
# 2) Re-colour the applied Design/theme from the "Integral" (Red Violet)
#    palette to the stock "Office Theme" palette.
# ---------------------------------------------------------------------------
function ToComRgb([string]$hex) {
    # PowerPoint's RGB() colour property packs bytes as R + G*256 + B*65536,
    # i.e. the reverse of how the hex string is written (RRGGBB).
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# Slot order exposed by ThemeColorScheme.Colors(1..12):
#   dk1, lt1, dk2, lt2, accent1, accent2, accent3, accent4, accent5,
#   accent6, hlink, folHlink
$officeThemeColors = "000000", "FFFFFF", "44546A", "E7E6E6",
                      "5B9BD5", "ED7D31", "A5A5A5", "FFC000",
                      "4472C4", "70AD47", "0563C1", "954F72"

$themeColors = $p.SlideMaster.Theme.ThemeColorScheme
for ($i = 1; $i -le 12; $i++) {
    $themeColors.Colors($i).RGB = ToComRgb($officeThemeColors[$i - 1])
}
